{"js": "const pairs = [\n  [\"2026-02-26 Thursday\", \"2026-02-27 Friday\"],\n  [\"78\u00d730=2340\", \"33\u00d741=1353\"],\n  [\"73\u00d753=3869\", \"64\u00d749=3136\"],\n  [\"99\u00d720=1980\", \"90\u00d799=8910\"],\n  [\"24\u00d799=2376\", \"58\u00d798=5684\"],\n  [\"25\u00d732=800\", \"24\u00d733=792\"],\n  [\"55\u00d755=3025\", \"71\u00d766=4686\"],\n  [\"46\u00d796=4416\", \"59\u00d789=5251\"],\n  [\"38\u00d769=2622\", \"23\u00d788=2024\"],\n  [\"88\u00d730=2640\", \"50\u00d761=3050\"],\n  [\"67\u00d724=1608\", \"32\u00d773=2336\"],\n  [\"47\u00d768=3196\", \"65\u00d729=1885\"],\n  [\"16\u00d779=1264\", \"14\u00d780=1120\"],\n  [\"60\u00d715=900\", \"53\u00d741=2173\"],\n  [\"66\u00d729=1914\", \"15\u00d728=420\"],\n  [\"28\u00d784=2352\", \"84\u00d743=3612\"],\n  [\"87\u00d729=2523\", \"34\u00d756=1904\"],\n  [\"36\u00d776=2736\", \"80\u00d761=4880\"],\n  [\"53\u00d786=4558\", \"77\u00d750=3850\"],\n  [\"24\u00d746=1104\", \"38\u00d714=532\"],\n  [\"18\u00d782=1476\", \"46\u00d745=2070\"],\n  [\"46\u00d768=3128\", \"53\u00d746=2438\"],\n  [\"65\u00d771=4615\", \"59\u00d736=2124\"],\n  [\"60\u00d717=1020\", \"94\u00d769=6486\"],\n  [\"90\u00d748=4320\", \"76\u00d727=2052\"],\n  [\"40\u00d753=2120\", \"18\u00d762=1116\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2026-02-26 Thursday', '2026-02-27 Friday'),\n    @('78\u00d730=2340', '33\u00d741=1353'),\n    @('73\u00d753=3869', '64\u00d749=3136'),\n    @('99\u00d720=1980', '90\u00d799=8910'),\n    @('24\u00d799=2376', '58\u00d798=5684'),\n    @('25\u00d732=800', '24\u00d733=792'),\n    @('55\u00d755=3025', '71\u00d766=4686'),\n    @('46\u00d796=4416', '59\u00d789=5251'),\n    @('38\u00d769=2622', '23\u00d788=2024'),\n    @('88\u00d730=2640', '50\u00d761=3050'),\n    @('67\u00d724=1608', '32\u00d773=2336'),\n    @('47\u00d768=3196', '65\u00d729=1885'),\n    @('16\u00d779=1264', '14\u00d780=1120'),\n    @('60\u00d715=900', '53\u00d741=2173'),\n    @('66\u00d729=1914', '15\u00d728=420'),\n    @('28\u00d784=2352', '84\u00d743=3612'),\n    @('87\u00d729=2523', '34\u00d756=1904'),\n    @('36\u00d776=2736', '80\u00d761=4880'),\n    @('53\u00d786=4558', '77\u00d750=3850'),\n    @('24\u00d746=1104', '38\u00d714=532'),\n    @('18\u00d782=1476', '46\u00d745=2070'),\n    @('46\u00d768=3128', '53\u00d746=2438'),\n    @('65\u00d771=4615', '59\u00d736=2124'),\n    @('60\u00d717=1020', '94\u00d769=6486'),\n    @('90\u00d748=4320', '76\u00d727=2052'),\n    @('40\u00d753=2120', '18\u00d762=1116'),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
